$wb = $excel.ActiveWorkbook

# Updates for sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 676
$ws1.Range("F4").Value = 38
$ws1.Range("F6").Value = 50
$ws1.Range("F7").Value = 46
$ws1.Range("F8").Value = 3225
$ws1.Range("F9").Value = 4217
$ws1.Range("F10").Value = 108

# Updates for sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 676
$ws4.Range("F4").Value = 38
$ws4.Range("F6").Value = 50
$ws4.Range("F7").Value = 46
$ws4.Range("F8").Value = 3225
$ws4.Range("F9").Value = 4217
$ws4.Range("F10").Value = 108
